$p = $ppt.ActivePresentation

$oldText = "https://www.youtube.com/watch?v=G3PvTWRIhZA&list=PLQVvvaa0QuDeF3hP0wQoSxpkqgRcgxMqX&index=1"
$newText = "ttps://www.youtube.com/watch?v=nSYFfWijl8U&index=2&list=PLQVvvaa0QuDeF3hP0wQoSxpkqgRcgxMqX"

$slideIndexes = @(2, 3, 4, 5)

foreach ($idx in $slideIndexes) {
    $s = $p.Slides.Item($idx)
    foreach ($shape in $s.Shapes) {
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}
